$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output Data")

# Row 2 (Atlantis)
$ws.Range("C2").Value = 700   # undiscounted reward
$ws.Range("L2").Value = 0     # country prob (R22)
$ws.Range("R2").Value = 0     # discounted reward (R21X)
$ws.Range("S2").Value = 0     # schedule prob (R22X)

# Row 4 (Carpania)
$ws.Range("L4").Value = 0     # expected utility (R22)
